$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D7").Value = 2
